# Swap the presentation's applied theme palette from the custom
# "Integral" (Red Violet) colour scheme over to the stock Office
# Theme colour scheme, via the Theme / ThemeColorScheme object model
# (ppt/theme/theme1.xml, the theme used by the slide master & all
# slides/layouts).
#
# PowerPoint's ThemeColorScheme index order (MsoThemeColorSchemeIndex):
#   1 = Dark1 (tx1)      2 = Light1 (bg1)
#   3 = Dark2 (tx2)      4 = Light2 (bg2)
#   5..10 = Accent1..6
#   11 = Hyperlink       12 = FollowedHyperlink

$p = $ppt.ActivePresentation

function RGBInt([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$officeColors = @(
    "000000",  # 1  dk1 / tx1
    "FFFFFF",  # 2  lt1 / bg1
    "44546A",  # 3  dk2 / tx2
    "E7E6E6",  # 4  lt2 / bg2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme

for ($i = 1; $i -le $officeColors.Length; $i++) {
    $themeColors.Item($i).RGB = RGBInt($officeColors[$i - 1])
}
